# Refactorizacion de estructura de scrapping
# Shift the weather forecast table forward by one week (dates +6 days,
# weekday labels shifted down one row) and update the temperature values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - Day labels
$ws.Range("A4").Value = "Sábado"
$ws.Range("A5").Value = "Domingo"
$ws.Range("A6").Value = "Lunes"
$ws.Range("A7").Value = "Martes"
$ws.Range("A8").Value = "Miércoles"

# Column B - Dates
$ws.Range("B2").Value = "8 Ago"
$ws.Range("B3").Value = "9 Ago"
$ws.Range("B4").Value = "10 Ago"
$ws.Range("B5").Value = "11 Ago"
$ws.Range("B6").Value = "12 Ago"
$ws.Range("B7").Value = "13 Ago"
$ws.Range("B8").Value = "14 Ago"

# Column C - Max/Min temperatures
$ws.Range("C2").Value = "38° / 23°"
$ws.Range("C3").Value = "36° / 23°"
$ws.Range("C4").Value = "36° / 23°"
$ws.Range("C5").Value = "37° / 22°"
$ws.Range("C6").Value = "37° / 21°"
$ws.Range("C7").Value = "37° / 22°"
$ws.Range("C8").Value = "37° / 23°"
